# Updates cryptos list values per commit: "Updated cryptos list on Fri May 19 06:12:59 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.914.38"
$ws.Range("E2").Value = "  -1.24%  "

# Row 3
$ws.Range("D3").Value = "1.807.84"
$ws.Range("E3").Value = "  -0.72%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").Value = "'310.19"
$ws.Range("E5").Value = "  -1.10%  "

# Row 6
$ws.Range("E6").Value = "  -0.08%  "

# Row 7
$ws.Range("D7").Value = "'0.4637"
$ws.Range("E7").Value = "  +3.22%  "

# Row 8
$ws.Range("E8").Value = "  -1.86%  "

# Row 9
$ws.Range("D9").Value = "'0.07383"
$ws.Range("E9").Value = "  -0.19%  "

# Row 10
$ws.Range("D10").Value = "'0.8746"
$ws.Range("E10").Value = "  -0.66%  "

# Row 11
$ws.Range("D11").Value = "'20.44"
$ws.Range("E11").Value = "  -1.93%  "

# Row 12
$ws.Range("D12").Value = "1.831.28"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13
$ws.Range("D13").Value = "'5.378"
$ws.Range("E13").Value = "  -0.96%  "

# Row 14
$ws.Range("D14").Value = "'92.65"
$ws.Range("E14").Value = "  -0.35%  "

# Row 15
$ws.Range("D15").Value = "'6.511"
$ws.Range("E15").Value = "  -2.94%  "

# Row 16
$ws.Range("D16").Value = "'0.07033"
$ws.Range("E16").Value = "  -1.36%  "

# Row 17
$ws.Range("E17").Value = "  -0.17%  "

# Row 18
$ws.Range("D18").Value = "'0.000008720"
$ws.Range("E18").Value = "  -0.87%  "

# Row 19
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.07%  "

# Row 20
$ws.Range("D20").Value = "'14.70"
$ws.Range("E20").Value = "  -2.42%  "

# Row 21
$ws.Range("D21").Value = "26.909.88"
$ws.Range("E21").Value = "  -1.30%  "

# Row 22
$ws.Range("D22").Value = "'5.309"
$ws.Range("E22").Value = "  -0.95%  "

# Row 23
$ws.Range("E23").Value = "  -2.52%  "

# Row 24
$ws.Range("D24").Value = "2.000.60"
$ws.Range("E24").Value = "  -2.41%  "

# Row 25
$ws.Range("D25").Value = "'1.906"
$ws.Range("E25").Value = "  -2.87%  "

# Row 26
$ws.Range("D26").Value = "'151.54"
$ws.Range("E26").Value = "  +0.07%  "

# Row 27
$ws.Range("D27").Value = "'18.35"
$ws.Range("E27").Value = "  -1.49%  "

# Row 28
$ws.Range("D28").Value = "'2.142"
$ws.Range("E28").Value = "  -6.31%  "

# Row 29
$ws.Range("D29").Value = "'5.304"
$ws.Range("E29").Value = "  -0.67%  "

# Row 30
$ws.Range("D30").Value = "'115.99"
$ws.Range("E30").Value = "  -1.08%  "

# Row 31
$ws.Range("D31").Value = "'0.08928"
$ws.Range("E31").Value = "  +0.47%  "

# Row 32
$ws.Range("D32").Value = "'0.7567"
$ws.Range("E32").Value = "  -3.26%  "

# Row 33
$ws.Range("D33").Value = "'1.158"
$ws.Range("E33").Value = "  -3.06%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.928"
$ws.Range("E34").Value = "  +0.78%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'4.463"
$ws.Range("E35").Value = "  -2.49%  "

# Row 36
$ws.Range("D36").Value = "'0.9999"
$ws.Range("E36").Value = "  -0.05%  "

# Row 37
$ws.Range("D37").Value = "'1.106"
$ws.Range("E37").Value = "  -0.44%  "

# Row 38
$ws.Range("D38").Value = "'0.01964"
$ws.Range("E38").Value = "  -0.60%  "

# Row 39
$ws.Range("D39").Value = "'0.05260"
$ws.Range("E39").Value = "  -0.29%  "

# Row 40
$ws.Range("D40").Value = "'2.407"
$ws.Range("E40").Value = "  +5.38%  "

# Row 41
$ws.Range("E41").Value = "  +1.95%  "

# Row 42
$ws.Range("D42").Value = "'7.221"
$ws.Range("E42").Value = "  -1.17%  "

# Row 43
$ws.Range("D43").Value = "'0.5303"
$ws.Range("E43").Value = "  +0.20%  "

# Row 44
$ws.Range("D44").Value = "'0.1670"
$ws.Range("E44").Value = "  -2.08%  "

# Row 45
$ws.Range("D45").Value = "'8.532"
$ws.Range("E45").Value = "  -0.71%  "

# Row 46
$ws.Range("E46").Value = "  -0.96%  "

# Row 47
$ws.Range("D47").Value = "'10.34"
$ws.Range("E47").Value = "  -2.09%  "

# Row 48
$ws.Range("D48").Value = "'103.92"
$ws.Range("E48").Value = "  -0.82%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.671"
$ws.Range("E49").Value = "  -0.85%  "

# Row 50
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'0.9998"
$ws.Range("E50").Value = "  -0.05%  "

# Row 51
$ws.Range("D51").Value = "'0.06293"
$ws.Range("E51").Value = "  -1.67%  "
